$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.803092333333333
$ws.Range("H2").Value = 5.409276999999999
$ws.Range("I2").Value = 0.1744886524959502
$ws.Range("J2").Value = 0.1744886524959502
$ws.Range("M2").Value = 3.078094333333333
$ws.Range("N2").Value = 9.234283
$ws.Range("O2").Value = 0.1162262227649079
$ws.Range("P2").Value = 0.1162262227649079
$ws.Range("Q2").Value = 5.55008829371011
$ws.Range("R2").Value = 49.95079464339099
$ws.Range("S2").Value = 0.02028015699494291
$ws.Range("T2").Value = 0.02028015699494291
$ws.Range("G3").Value = 1.803092333333333
$ws.Range("H3").Value = 5.409276999999999
$ws.Range("I3").Value = 0.1744886524959502
$ws.Range("J3").Value = 0.1744886524959502
$ws.Range("O3").Value = 0.390990551238643
$ws.Range("P3").Value = 0.390990551238643
$ws.Range("Q3").Value = 18.67076146637066
$ws.Range("R3").Value = 168.036853197336
$ws.Range("S3").Value = 0.06822341442427961
$ws.Range("T3").Value = 0.06822341442427959
$ws.Range("G4").Value = 1.803092333333333
$ws.Range("H4").Value = 5.409276999999999
$ws.Range("I4").Value = 0.1744886524959502
$ws.Range("J4").Value = 0.1744886524959502
$ws.Range("M4").Value = 12.706793
$ws.Range("N4").Value = 38.120379
$ws.Range("O4").Value = 0.4797976910104138
$ws.Range("P4").Value = 0.4797976910104139
$ws.Range("Q4").Value = 22.91152103955366
$ws.Range("R4").Value = 206.203689355983
$ws.Range("S4").Value = 0.08371925257507541
$ws.Range("T4").Value = 0.0837192525750754
$ws.Range("G5").Value = 1.803092333333333
$ws.Range("H5").Value = 5.409276999999999
$ws.Range("I5").Value = 0.1744886524959502
$ws.Range("J5").Value = 0.1744886524959502
$ws.Range("M5").Value = 0.3439043333333334
$ws.Range("N5").Value = 1.031713
$ws.Range("O5").Value = 0.01298553498603535
$ws.Range("P5").Value = 0.01298553498603535
$ws.Range("Q5").Value = 0.6200912668334444
$ws.Range("R5").Value = 5.580821401501
$ws.Range("S5").Value = 0.002265828501652326
$ws.Range("T5").Value = 0.002265828501652325
$ws.Range("I6").Value = 0.4384883998568034
$ws.Range("J6").Value = 0.4384883998568034
$ws.Range("M6").Value = 3.078094333333333
$ws.Range("N6").Value = 9.234283
$ws.Range("O6").Value = 0.1162262227649079
$ws.Range("P6").Value = 0.1162262227649079
$ws.Range("Q6").Value = 13.94732150292355
$ws.Range("R6").Value = 125.525893526312
$ws.Range("S6").Value = 0.05096385044158484
$ws.Range("T6").Value = 0.05096385044158484
$ws.Range("I7").Value = 0.4384883998568034
$ws.Range("J7").Value = 0.4384883998568034
$ws.Range("O7").Value = 0.390990551238643
$ws.Range("P7").Value = 0.390990551238643
$ws.Range("R7").Value = 422.275086783552
$ws.Range("S7").Value = 0.1714448211717621
$ws.Range("T7").Value = 0.171444821171762
$ws.Range("I8").Value = 0.4384883998568034
$ws.Range("J8").Value = 0.4384883998568034
$ws.Range("M8").Value = 12.706793
$ws.Range("N8").Value = 38.120379
$ws.Range("O8").Value = 0.4797976910104138
$ws.Range("P8").Value = 0.4797976910104139
$ws.Range("Q8").Value = 57.57644440031733
$ws.Range("R8").Value = 518.1879996028559
$ws.Range("S8").Value = 0.2103857217861453
$ws.Range("T8").Value = 0.2103857217861453
$ws.Range("I9").Value = 0.4384883998568034
$ws.Range("J9").Value = 0.4384883998568034
$ws.Range("M9").Value = 0.3439043333333334
$ws.Range("N9").Value = 1.031713
$ws.Range("O9").Value = 0.01298553498603535
$ws.Range("P9").Value = 0.01298553498603535
$ws.Range("Q9").Value = 1.558283724870222
$ws.Range("R9").Value = 14.024553523832
$ws.Range("S9").Value = 0.005694006457311177
$ws.Range("T9").Value = 0.005694006457311177
$ws.Range("G10").Value = 3.895605666666667
$ws.Range("H10").Value = 11.686817
$ws.Range("I10").Value = 0.3769851220961256
$ws.Range("J10").Value = 0.3769851220961256
$ws.Range("M10").Value = 3.078094333333333
$ws.Range("N10").Value = 9.234283
$ws.Range("O10").Value = 0.1162262227649079
$ws.Range("P10").Value = 0.1162262227649079
$ws.Range("Q10").Value = 11.99104172746789
$ws.Range("R10").Value = 107.919375547211
$ws.Range("S10").Value = 0.04381555677980029
$ws.Range("T10").Value = 0.0438155567798003
$ws.Range("G11").Value = 3.895605666666667
$ws.Range("H11").Value = 11.686817
$ws.Range("I11").Value = 0.3769851220961256
$ws.Range("J11").Value = 0.3769851220961256
$ws.Range("O11").Value = 0.390990551238643
$ws.Range("P11").Value = 0.390990551238643
$ws.Range("Q11").Value = 40.33843571111733
$ws.Range("R11").Value = 363.045921400056
$ws.Range("S11").Value = 0.1473976206971313
$ws.Range("T11").Value = 0.1473976206971313
$ws.Range("G12").Value = 3.895605666666667
$ws.Range("H12").Value = 11.686817
$ws.Range("I12").Value = 0.3769851220961256
$ws.Range("J12").Value = 0.3769851220961256
$ws.Range("M12").Value = 12.706793
$ws.Range("N12").Value = 38.120379
$ws.Range("O12").Value = 0.4797976910104138
$ws.Range("P12").Value = 0.4797976910104139
$ws.Range("Q12").Value = 49.50065481596033
$ws.Range("R12").Value = 445.505893343643
$ws.Range("S12").Value = 0.180876591127
$ws.Range("T12").Value = 0.180876591127
$ws.Range("G13").Value = 3.895605666666667
$ws.Range("H13").Value = 11.686817
$ws.Range("I13").Value = 0.3769851220961256
$ws.Range("J13").Value = 0.3769851220961256
$ws.Range("M13").Value = 0.3439043333333334
$ws.Range("N13").Value = 1.031713
$ws.Range("O13").Value = 0.01298553498603535
$ws.Range("P13").Value = 0.01298553498603535
$ws.Range("Q13").Value = 1.339715669724556
$ws.Range("R13").Value = 12.057441027521
$ws.Range("S13").Value = 0.004895353492194046
$ws.Range("T13").Value = 0.004895353492194046
$ws.Range("G14").Value = 0.1037266666666667
$ws.Range("H14").Value = 0.31118
$ws.Range("I14").Value = 0.01003782555112075
$ws.Range("J14").Value = 0.01003782555112075
$ws.Range("M14").Value = 3.078094333333333
$ws.Range("N14").Value = 9.234283
$ws.Range("O14").Value = 0.1162262227649079
$ws.Range("P14").Value = 0.1162262227649079
$ws.Range("Q14").Value = 0.3192804648822222
$ws.Range("R14").Value = 2.87352418394
$ws.Range("S14").Value = 0.001166658548579845
$ws.Range("T14").Value = 0.001166658548579845
$ws.Range("G15").Value = 0.1037266666666667
$ws.Range("H15").Value = 0.31118
$ws.Range("I15").Value = 0.01003782555112075
$ws.Range("J15").Value = 0.01003782555112075
$ws.Range("O15").Value = 0.390990551238643
$ws.Range("P15").Value = 0.390990551238643
$ws.Range("Q15").Value = 1.074074696693333
$ws.Range("R15").Value = 9.666672270240001
$ws.Range("S15").Value = 0.003924694945470038
$ws.Range("T15").Value = 0.003924694945470037
$ws.Range("G16").Value = 0.1037266666666667
$ws.Range("H16").Value = 0.31118
$ws.Range("I16").Value = 0.01003782555112075
$ws.Range("J16").Value = 0.01003782555112075
$ws.Range("M16").Value = 12.706793
$ws.Range("N16").Value = 38.120379
$ws.Range("O16").Value = 0.4797976910104138
$ws.Range("P16").Value = 0.4797976910104139
$ws.Range("Q16").Value = 1.318033281913333
$ws.Range("R16").Value = 11.86229953722
$ws.Range("S16").Value = 0.004816125522193071
$ws.Range("T16").Value = 0.004816125522193071
$ws.Range("G17").Value = 0.1037266666666667
$ws.Range("H17").Value = 0.31118
$ws.Range("I17").Value = 0.01003782555112075
$ws.Range("J17").Value = 0.01003782555112075
$ws.Range("M17").Value = 0.3439043333333334
$ws.Range("N17").Value = 1.031713
$ws.Range("O17").Value = 0.01298553498603535
$ws.Range("P17").Value = 0.01298553498603535
$ws.Range("Q17").Value = 0.0356720501488889
$ws.Range("R17").Value = 0.32104845134
$ws.Range("S17").Value = 0.000130346534877798
$ws.Range("T17").Value = 0.0001316663680575393
